# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the existing summary column (G) and filling in the
# per-row save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: H1 = "Save", formatted like the other header cells (G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows: H2:H9 hold the save indicator for each row.
$saveValues = @(0, 0, 0, 1, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
